$d = $word.ActiveDocument

# --- N02: simple text swap (keeps existing run/paragraph formatting) ---
$d.Content.Find.Execute(
    "N02: Monitoramento contínuo do estoque;", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "N02: Venda de forma rápida e eficaz através da loja virtual;", 2)

# --- N03: simple text swap (keeps existing run/paragraph formatting) ---
$d.Content.Find.Execute(
    "N03: Comunicação com os fornecedores de forma automatizada;", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "N03: Segurança;", 2)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- N04 (was "Venda de forma rápida..."): new text, drop the paragraph-mark
#     underline override (<w:u w:val="none"/>) that lived in pPr/rPr ---
$p4 = $d.Paragraphs.Item(6)
$p4.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">N04: Controle do catálogo online;</w:t></w:r></w:p>
"@)

# --- N05 (was "Suporte a vendas de alta demanda;"): new text + trailing ';',
#     drop the paragraph-mark underline override ---
$p5 = $d.Paragraphs.Item(7)
$p5.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">N05: Cadastro de cliente;</w:t></w:r></w:p>
"@)

# --- N06 (was "Controle do catálogo online;"): new text, drop the
#     paragraph-mark underline override ---
$p6 = $d.Paragraphs.Item(8)
$p6.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">N06: Suporte de vendas.</w:t></w:r></w:p>
"@)

# --- N07 (was "Cadastro de cliente"): new text, drop the paragraph-mark
#     underline override ---
$p7 = $d.Paragraphs.Item(9)
$p7.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">N07: Administração de empresa</w:t></w:r></w:p>
"@)

# --- former N08 ("Rastreamento de entrega.") paragraph becomes a blank,
#     non-list paragraph (left indent kept, no hanging/numbering) ---
$p8 = $d.Paragraphs.Item(10)
$p8.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:ind w:left="720" w:firstLine="0"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>
"@)
# force the explicit firstLine="0" (instead of just "no first-line indent")
# to stick on the paragraph mark, matching the target markup exactly
$d.Paragraphs.Item(10).Format.FirstLineIndent = 0
